# Apply updated crypto price/volume data to match target snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "35.426.99"
$ws.Cells.Item(2, 5).Value = "  +0.41%  "

$ws.Cells.Item(3, 4).Value = "1.884.54"
$ws.Cells.Item(3, 5).Value = "  -1.49%  "

$ws.Cells.Item(4, 5).Value = "  -0.93%  "

$ws.Cells.Item(5, 4).Value = "'246.02"
$ws.Cells.Item(5, 5).Value = "  -3.83%  "

$ws.Cells.Item(6, 4).Value = "'0.689"
$ws.Cells.Item(6, 5).Value = "  -4.93%  "

$ws.Cells.Item(7, 5).Value = "  -0.87%  "

$ws.Cells.Item(8, 4).Value = "'43.13"
$ws.Cells.Item(8, 5).Value = "  +2.81%  "

$ws.Cells.Item(9, 5).Value = "  -4.30%  "

$ws.Cells.Item(10, 5).Value = "  -3.49%  "

$ws.Cells.Item(11, 4).Value = "'0.0970"
$ws.Cells.Item(11, 5).Value = "  -1.94%  "

$ws.Cells.Item(12, 4).Value = "'13.08"
$ws.Cells.Item(12, 5).Value = "  -0.40%  "

$ws.Cells.Item(13, 4).Value = "2.156.78"
$ws.Cells.Item(13, 5).Value = "  -1.70%  "

$ws.Cells.Item(14, 4).Value = "'0.741"
$ws.Cells.Item(14, 5).Value = "  +0.54%  "

$ws.Cells.Item(15, 5).Value = "  -1.34%  "

$ws.Cells.Item(16, 4).Value = "1.884.12"
$ws.Cells.Item(16, 5).Value = "  -1.69%  "

$ws.Cells.Item(17, 4).Value = "35.356.92"
$ws.Cells.Item(17, 5).Value = "  +0.14%  "

$ws.Cells.Item(18, 4).Value = "'73.49"
$ws.Cells.Item(18, 5).Value = "  -2.32%  "

$ws.Cells.Item(19, 4).Value = "0.0₃0822"
$ws.Cells.Item(19, 5).Value = "  -3.21%  "

$ws.Cells.Item(20, 4).Value = "'245.16"
$ws.Cells.Item(20, 5).Value = "  -0.44%  "

$ws.Cells.Item(21, 5).Value = "  -2.29%  "

$ws.Cells.Item(22, 4).Value = "'4.92"
$ws.Cells.Item(22, 5).Value = "  -4.49%  "

$ws.Cells.Item(23, 5).Value = "  -0.72%  "

$ws.Cells.Item(24, 4).Value = "'2.54"
$ws.Cells.Item(24, 5).Value = "  +3.50%  "

$ws.Cells.Item(25, 5).Value = "  -9.87%  "

$ws.Cells.Item(26, 4).Value = "'165.61"
$ws.Cells.Item(26, 5).Value = "  -0.51%  "

$ws.Cells.Item(27, 4).Value = "'8.47"
$ws.Cells.Item(27, 5).Value = "  -3.71%  "

$ws.Cells.Item(28, 4).Value = "'18.33"
$ws.Cells.Item(28, 5).Value = "  -2.52%  "

$ws.Cells.Item(29, 5).Value = "  -4.20%  "

$ws.Cells.Item(30, 4).Value = "4.128.42"
$ws.Cells.Item(30, 5).Value = "  -0.01%  "

$ws.Cells.Item(31, 4).Value = "'1.74"
$ws.Cells.Item(31, 5).Value = "  +3.56%  "

$ws.Cells.Item(32, 5).Value = "  -3.16%  "

$ws.Cells.Item(33, 4).Value = "'0.0578"
$ws.Cells.Item(33, 5).Value = "  -2.40%  "

$ws.Cells.Item(34, 4).Value = "'4.20"
$ws.Cells.Item(34, 5).Value = "  -1.79%  "

$ws.Cells.Item(35, 5).Value = "  -1.00%  "

$ws.Cells.Item(36, 4).Value = "'0.850"
$ws.Cells.Item(36, 5).Value = "  -7.29%  "

$ws.Cells.Item(37, 4).Value = "'1.66"
$ws.Cells.Item(37, 5).Value = "  -17.09%  "

$ws.Cells.Item(38, 5).Value = "  -3.76%  "

$ws.Cells.Item(39, 4).Value = "'0.0683"
$ws.Cells.Item(39, 5).Value = "  +4.83%  "

$ws.Cells.Item(40, 4).Value = "'97.16"
$ws.Cells.Item(40, 5).Value = "  -2.48%  "

$ws.Cells.Item(41, 2).Value = "VeChain"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(41, 4).Value = "'0.0216"
$ws.Cells.Item(41, 5).Value = "  -2.23%  "

$ws.Cells.Item(42, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(42, 4).Value = "'16.98"
$ws.Cells.Item(42, 5).Value = "  -0.32%  "

$ws.Cells.Item(43, 5).Value = "  -4.17%  "

$ws.Cells.Item(44, 4).Value = "1.293.31"
$ws.Cells.Item(44, 5).Value = "  -3.81%  "

$ws.Cells.Item(45, 4).Value = "'2.34"
$ws.Cells.Item(45, 5).Value = "  -4.89%  "

$ws.Cells.Item(46, 5).Value = "  +6.78%  "

$ws.Cells.Item(47, 5).Value = "  -1.29%  "

$ws.Cells.Item(48, 5).Value = "  -0.91%  "

$ws.Cells.Item(49, 4).Value = "'12.13"
$ws.Cells.Item(49, 5).Value = "  +3.88%  "

$ws.Cells.Item(50, 4).Value = "'43.03"
$ws.Cells.Item(50, 5).Value = "  -4.46%  "

$ws.Cells.Item(51, 4).Value = "'6.22"
$ws.Cells.Item(51, 5).Value = "  -7.78%  "
